$wb = $excel.ActiveWorkbook

# --- Sheet at position 1 currently holds "hotel_info" data; it will become
#     the "review_info" sheet (headers only, no data rows).
# --- Sheet at position 2 currently holds "review_info" headers; it will
#     become the "hotel_info" sheet (headers + one data row, with a new
#     "State" column inserted after "Hotel_Name").
#
# We swap the *content* of the two physical sheets (rather than moving the
# tabs) so that rId1/sheetId1 ends up named "review_info" and
# rId2/sheetId2 ends up named "hotel_info", matching tab order.

$wsA = $wb.Worksheets.Item(1)
$wsB = $wb.Worksheets.Item(2)

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

# Row 2 values for the hotel_info sheet. $null marks columns written as
# real numbers; the rest are written as text (matching the source data,
# where counters like "771" are stored as text, not numbers).
$hotelRow2 = @(
    40374,
    "Hilton Garden Inn New Orleans Airport",
    "Louisiana",
    "Kenner",
    70065,
    "https://www.tripadvisor.com/Hotel_Review-g40247-d224644-Reviews-Hilton_Garden_Inn_New_Orleans_Airport-Kenner_Louisiana.html",
    "Hilton Garden Inn New Orleans Airport",
    "771",
    "6",
    "774"
)

# Columns (1-based) in $hotelRow2 that must be stored as text even though
# they look numeric.
$hotelRow2TextCols = @(8, 9, 10)

# Rename to temporary, collision-free names first (both current names get
# reused, just swapped between the two physical sheets).
$wsA.Name = "__tmp_sheet_A__"
$wsB.Name = "__tmp_sheet_B__"

# --- Rebuild sheet at position 1 as "review_info" (headers only) ----------
$wsA.Cells.Clear()
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $wsA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}
$wsA.Name = "review_info"

# --- Rebuild sheet at position 2 as "hotel_info" (headers + data row) -----
$wsB.Cells.Clear()
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $wsB.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}
for ($i = 0; $i -lt $hotelRow2.Length; $i++) {
    $col = $i + 1
    $cell = $wsB.Cells.Item(2, $col)
    if ($hotelRow2TextCols -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $hotelRow2[$i]
}
$wsB.Name = "hotel_info"
